# Auto-generated edit script applying profit/price recalculations
# from the scheduled market-data runner update.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 605.55554
$ws.Range("I28").Value = 431.2143
$ws.Range("K28").Value = 431.2143
$ws.Range("M28").Value = 53.78570000000002
$ws.Range("H106").Value = 2922.5
$ws.Range("I106").Value = 2896.6667
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 2896.6667
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -2265.6667
$ws.Range("N106").Value = -4262
$ws.Range("H107").Value = 1316.4762
$ws.Range("I107").Value = 1112.5
$ws.Range("K107").Value = 1112.5
$ws.Range("M107").Value = 807.5
$ws.Range("H115").Value = 368.8889
$ws.Range("I115").Value = 368.8889
$ws.Range("K115").Value = 1106.6667
$ws.Range("M115").Value = 460.3333
$ws.Range("H138").Value = 2238.5103
$ws.Range("I138").Value = 843.7273
$ws.Range("J138").Value = 3375
$ws.Range("K138").Value = 2531.1819
$ws.Range("L138").Value = 10125
$ws.Range("M138").Value = 2608.8181
$ws.Range("N138").Value = -20405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3502.8215
$ws.Range("I32").Value = 3655.157
$ws.Range("K32").Value = 3655.157
$ws.Range("M32").Value = -3368.157
$ws.Range("H45").Value = 4185.125
$ws.Range("I45").Value = 3538.6667
$ws.Range("K45").Value = 3538.6667
$ws.Range("M45").Value = -3161.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 28727.092
$ws.Range("I86").Value = 12666.444
$ws.Range("K86").Value = 12666.444
$ws.Range("M86").Value = -11543.444
$ws.Range("H89").Value = 28727.092
$ws.Range("I89").Value = 12666.444
$ws.Range("K89").Value = 63332.22
$ws.Range("M89").Value = -57716.22
$ws.Range("H107").Value = 1813.1515
$ws.Range("I107").Value = 1672.32
$ws.Range("J107").Value = 2253.25
$ws.Range("K107").Value = 1672.32
$ws.Range("L107").Value = 2253.25
$ws.Range("M107").Value = 247.6800000000001
$ws.Range("N107").Value = -6093.25
$ws.Range("H134").Value = 2732.7368
$ws.Range("I134").Value = 1806.8125
$ws.Range("K134").Value = 5420.4375
$ws.Range("M134").Value = -2885.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 112.07692
$ws.Range("J7").Value = 158.83333
$ws.Range("L7").Value = 158.83333
$ws.Range("N7").Value = -384.83333
$ws.Range("H31").Value = 3431.7
$ws.Range("I31").Value = 1773
$ws.Range("K31").Value = 1773
$ws.Range("M31").Value = -1478
$ws.Range("H34").Value = 3431.7
$ws.Range("I34").Value = 1773
$ws.Range("K34").Value = 1773
$ws.Range("M34").Value = -1571
$ws.Range("H94").Value = 2830.1667
$ws.Range("I94").Value = 2999
$ws.Range("J94").Value = 2814.818
$ws.Range("K94").Value = 2999
$ws.Range("L94").Value = 2814.818
$ws.Range("M94").Value = -2548
$ws.Range("N94").Value = -3716.818
$ws.Range("H107").Value = 756
$ws.Range("J107").Value = 756
$ws.Range("L107").Value = 756
$ws.Range("N107").Value = -4596
$ws.Range("H132").Value = 3190.975
$ws.Range("I132").Value = 3279.9143
$ws.Range("K132").Value = 9839.742899999999
$ws.Range("M132").Value = -7309.742899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3484
$ws.Range("I3").Value = 1980.125
$ws.Range("K3").Value = 5940.375
$ws.Range("M3").Value = -5828.375
$ws.Range("H5").Value = 1929.7778
$ws.Range("I5").Value = 195.42857
$ws.Range("K5").Value = 586.28571
$ws.Range("M5").Value = -474.28571
$ws.Range("H36").Value = 6670114.5
$ws.Range("I36").Value = 589.5
$ws.Range("J36").Value = 8337495.5
$ws.Range("K36").Value = 1768.5
$ws.Range("L36").Value = 25012486.5
$ws.Range("M36").Value = -1599.5
$ws.Range("N36").Value = -25012824.5
$ws.Range("H64").Value = 21600
$ws.Range("H67").Value = 21600
$ws.Range("H131").Value = 19232074
$ws.Range("I131").Value = 38462400
$ws.Range("K131").Value = 115387200
$ws.Range("M131").Value = -115382160
$ws.Range("H133").Value = 2815.6667
$ws.Range("I133").Value = 1943
$ws.Range("K133").Value = 5829
$ws.Range("M133").Value = -769
$ws.Range("H134").Value = 1280.48
$ws.Range("J134").Value = 2852
$ws.Range("L134").Value = 8556
$ws.Range("N134").Value = -18696
$ws.Range("H135").Value = 1929.7778
$ws.Range("I135").Value = 195.42857
$ws.Range("K135").Value = 1758.85713
$ws.Range("M135").Value = 776.1428699999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2056.9
$ws.Range("I97").Value = 931.6
$ws.Range("J97").Value = 3182.2
$ws.Range("K97").Value = 931.6
$ws.Range("L97").Value = 3182.2
$ws.Range("M97").Value = -435.6
$ws.Range("N97").Value = -4174.2
$ws.Range("H132").Value = 2192.3
$ws.Range("I132").Value = 1599.68
$ws.Range("J132").Value = 5155.4
$ws.Range("K132").Value = 4799.04
$ws.Range("L132").Value = 15466.2
$ws.Range("M132").Value = -2269.04
$ws.Range("N132").Value = -20526.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4199.75
$ws.Range("I93").Value = 4000
$ws.Range("J93").Value = 4799
$ws.Range("K93").Value = 4000
$ws.Range("L93").Value = 4799
$ws.Range("M93").Value = -2752
$ws.Range("N93").Value = -7295
$ws.Range("H122").Value = 4284.0293
$ws.Range("I122").Value = 4014.56
$ws.Range("J122").Value = 5032.5557
$ws.Range("K122").Value = 12043.68
$ws.Range("L122").Value = 15097.6671
$ws.Range("M122").Value = -9593.68
$ws.Range("N122").Value = -19997.6671
$ws.Range("H136").Value = 6137.1113
$ws.Range("I136").Value = 3520.6924
$ws.Range("J136").Value = 12939.8
$ws.Range("K136").Value = 10562.0772
$ws.Range("L136").Value = 38819.39999999999
$ws.Range("M136").Value = -8012.0772
$ws.Range("N136").Value = -43919.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 50000
$ws.Range("J27").Value = 50000
$ws.Range("L27").Value = 50000
$ws.Range("N27").Value = -50138
$ws.Range("H40").Value = 11655.8
$ws.Range("I40").Value = 7759.6665
$ws.Range("J40").Value = 17500
$ws.Range("K40").Value = 7759.6665
$ws.Range("L40").Value = 17500
$ws.Range("M40").Value = -7610.6665
$ws.Range("N40").Value = -17798
$ws.Range("H81").Value = 5009.375
$ws.Range("I81").Value = 3149.3635
$ws.Range("J81").Value = 7282.722
$ws.Range("K81").Value = 6298.727
$ws.Range("L81").Value = 14565.444
$ws.Range("M81").Value = -5237.727
$ws.Range("N81").Value = -16687.444
$ws.Range("H84").Value = 5009.375
$ws.Range("I84").Value = 3149.3635
$ws.Range("J84").Value = 7282.722
$ws.Range("K84").Value = 31493.635
$ws.Range("L84").Value = 72827.22
$ws.Range("M84").Value = -26189.635
$ws.Range("N84").Value = -83435.22
$ws.Range("H122").Value = 4692.625
$ws.Range("I122").Value = 3367.8462
$ws.Range("K122").Value = 10103.5386
$ws.Range("M122").Value = -7653.5386
